$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-10-17 Tuesday" "2023-10-18 Wednesday"

Replace-Text "93×75=6975" "78×47=3666"
Replace-Text "38×39=1482" "44×29=1276"
Replace-Text "12×64=768" "77×67=5159"
Replace-Text "97×39=3783" "17×11=187"
Replace-Text "31×14=434" "14×90=1260"
Replace-Text "54×63=3402" "37×27=999"
Replace-Text "71×89=6319" "97×42=4074"
Replace-Text "29×40=1160" "98×69=6762"
Replace-Text "80×14=1120" "41×86=3526"
Replace-Text "89×64=5696" "58×13=754"
Replace-Text "31×83=2573" "28×85=2380"
Replace-Text "49×21=1029" "66×92=6072"
Replace-Text "35×55=1925" "32×84=2688"
Replace-Text "56×52=2912" "47×95=4465"
Replace-Text "11×60=660" "77×96=7392"
Replace-Text "47×88=4136" "86×67=5762"
Replace-Text "11×99=1089" "94×24=2256"
Replace-Text "64×70=4480" "47×25=1175"
Replace-Text "28×59=1652" "38×81=3078"
Replace-Text "73×84=6132" "32×59=1888"
Replace-Text "43×16=688" "75×47=3525"
Replace-Text "62×30=1860" "80×58=4640"
Replace-Text "96×17=1632" "61×88=5368"
Replace-Text "79×79=6241" "39×78=3042"
Replace-Text "50×33=1650" "56×68=3808"
